$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CancelacionAhorros")

# Update header D1 from "razon" to "razon cierre"
$ws.Range("D1").Value = "razon cierre"

# Update the active selection on the sheet (cosmetic, matches author's navigation)
$ws.Range("F7").Select()
